$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf("2018")
$chars = $tr.Characters($idx+1, 4)
$chars.Text = "2019"
